# Fill in the blank "Status" column (D) values for the charity list on Sheet1.
# These are the rows that previously had no Status set.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$statusByRow = @{
    137 = "Done"
    138 = "Done"
    139 = "Done"
    140 = "Done"
    141 = "No Events"
    142 = "No Events"
    143 = "Done"
    144 = "No Events"
    145 = "No Events"
    146 = "Done"
    147 = "Done"
    148 = "Pending"
    149 = "No Events"
    150 = "No Events"
    151 = "No Events"
    152 = "Website Not Found"
    153 = "Done"
    154 = "Done"
    155 = "No Events"
    156 = "No Events"
    157 = "Done"
    158 = "No Events"
    159 = "No Events"
    160 = "Done"
    161 = "No Events"
    162 = "No Events"
    163 = "No Events"
    164 = "No Events"
    165 = "Done"
    166 = "Done"
    167 = "No Events"
    168 = "No Events"
    169 = "Done"
    170 = "No Events"
    171 = "No Events"
    172 = "No Events"
    173 = "Done"
    174 = "Pending"
    185 = "No Events"
    186 = "Large Site Scrapper needed"
    187 = "Done"
}

foreach ($row in $statusByRow.Keys) {
    $ws.Cells.Item($row, 4).Value = $statusByRow[$row]
}
